$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric literals need to be forced to
# Text format first, otherwise Excel auto-converts the assigned string into
# a number (losing the inline-string / text representation used in the file).
$textForceCells = @('D5', 'D6', 'D8', 'D10', 'D11', 'D14', 'D16', 'D20', 'D24', 'D25', 'D29', 'D35', 'D39', 'D41', 'D43', 'D46', 'D51')
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = '@'
}

$ws.Range('D2').Value = '26.666.31'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').Value = '1.599.05'
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').Value = '211.75'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').Value = '0.516'
$ws.Range('E6').Value = '  +1.01%  '
$ws.Range('D8').Value = '0.0619'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('D10').Value = '19.58'
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('D11').Value = '0.0838'
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').Value = '1.823.38'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.588.43'
$ws.Range('E13').Value = '  +3.03%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '4.03'
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('D16').Value = '65.26'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').Value = '26.669.51'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = '0.0₃0739'
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('D20').Value = '209.28'
$ws.Range('E20').Value = '  -0.30%  '
$ws.Range('E21').Value = '  +4.71%  '
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('D24').Value = '8.98'
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('D25').Value = '145.36'
$ws.Range('E25').Value = '  -1.00%  '
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('E27').Value = '  -1.00%  '
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('D29').Value = '15.31'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').Value = '  +2.11%  '
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('E32').Value = '  +0.25%  '
$ws.Range('E33').Value = '  +1.32%  '
$ws.Range('D34').Value = '1.280.06'
$ws.Range('E34').Value = '  -1.55%  '
$ws.Range('D35').Value = '0.621'
$ws.Range('E35').Value = '  -7.82%  '
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('E37').Value = '  +0.98%  '
$ws.Range('E38').Value = '  -0.95%  '
$ws.Range('D39').Value = '0.836'
$ws.Range('E39').Value = '  -1.01%  '
$ws.Range('E40').Value = '  +19.45%  '
$ws.Range('D41').Value = '5.52'
$ws.Range('E41').Value = '  +2.70%  '
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('D43').Value = '0.785'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').Value = '1.735.97'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').Value = '91.53'
$ws.Range('E46').Value = '  +1.61%  '
$ws.Range('E47').Value = '  -2.42%  '
$ws.Range('E48').Value = '  +3.43%  '
$ws.Range('E49').Value = '  +0.61%  '
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('D51').Value = '7.42'
$ws.Range('E51').Value = '  -1.39%  '
